$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1719.057593478942
$ws.Range("B3").Value = 2586.576303619814
$ws.Range("B4").Value = 3280.31833847926
$ws.Range("B5").Value = 3561.767802847471
$ws.Range("B6").Value = 3765.741480360838
$ws.Range("B7").Value = 3869.290887425631
$ws.Range("B8").Value = 3910.131292998394
$ws.Range("B9").Value = 3883.850855019931
$ws.Range("B10").Value = 3795.55919237603
$ws.Range("B11").Value = 3724.603231125617
$ws.Range("B12").Value = 3566.001545407796
$ws.Range("B13").Value = 3342.819166823189
$ws.Range("B14").Value = 3107.375746679888
$ws.Range("B15").Value = 2915.014578490384
$ws.Range("B16").Value = 2667.40687120604
$ws.Range("B17").Value = 2230.226820097469
$ws.Range("B18").Value = 2043.779951158386
$ws.Range("B19").Value = 1670.228475302213
$ws.Range("B20").Value = 1268.781579200319
$ws.Range("B21").Value = 1129.452597178012
$ws.Range("B22").Value = 563.4946251306228
$ws.Range("B23").Value = 395.678556308223
$ws.Range("B24").Value = 278.2820281813218
$ws.Range("B25").Value = 19.78964520059605
$ws.Range("B26").Value = 19.78964506479923
$ws.Range("B27").Value = 19.78964500483934
$ws.Range("B28").Value = 19.78964500483934
$ws.Range("B29").Value = 19.78964500483934
$ws.Range("B30").Value = 19.78964500483934
$ws.Range("B31").Value = 19.78964500483934
$ws.Range("B32").Value = 19.78964500483934
$ws.Range("B33").Value = 19.78964500483934
$ws.Range("B34").Value = 19.78964500483934
$ws.Range("B35").Value = 19.78964500483934
$ws.Range("B36").Value = 19.78964500483934
$ws.Range("B37").Value = 19.78964500483934
$ws.Range("B38").Value = 19.78964500483934
$ws.Range("B39").Value = 19.78964500483934
$ws.Range("B40").Value = 19.78964500483934
$ws.Range("B41").Value = 19.78964500483934
$ws.Range("B42").Value = 19.78964500483934
$ws.Range("B43").Value = 19.78964500483934
$ws.Range("B44").Value = 19.78964500483934
$ws.Range("B45").Value = 19.78964500483934
$ws.Range("B46").Value = 19.78964500483934
$ws.Range("B47").Value = 19.78964500483934
$ws.Range("B48").Value = 19.78964500483934
$ws.Range("B49").Value = 19.78964500483934
$ws.Range("B50").Value = 19.78964500483934
$ws.Range("B51").Value = 19.78964500483934
$ws.Range("B52").Value = 19.78964500483934
$ws.Range("B53").Value = 19.78964500483934
$ws.Range("B54").Value = 19.78964500483934
$ws.Range("B55").Value = 19.78964500483934
$ws.Range("B56").Value = 19.78964500483934
$ws.Range("B57").Value = 19.78964500483934
$ws.Range("B58").Value = 19.78964500483934
$ws.Range("B59").Value = 19.78964500483934
$ws.Range("B60").Value = 19.78964500483934
$ws.Range("B61").Value = 19.78964500483934
$ws.Range("B62").Value = 19.78964500483934
